# Week 3 day 1 — append 5 new tag-log rows (Action # 61-64, plus a trailing
# un-numbered row) to the bottom of the log table and grow Table1 to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# --- Row 62 : Action 61 -----------------------------------------------
$ws.Range("A62").Value = 61
$ws.Range("B62").Value = "Tag "
$ws.Range("C62").Value = "Tenel Enns"
$ws.Range("D62").Value = "Eliot Rogers"
$ws.Range("E62").Value = 45723
$ws.Range("F62").Value = 0.5541666666666667

# --- Row 63 : Action 62 -----------------------------------------------
$ws.Range("A63").Value = 62
$ws.Range("B63").Value = "Tag"
$ws.Range("C63").Value = "Alex Wates"
$ws.Range("D63").Value = "Ethan Darby"
$ws.Range("E63").Value = 45723
$ws.Range("F63").Value = 0.5541666666666667

# --- Row 64 : Action 63 -----------------------------------------------
$ws.Range("A64").Value = 63
$ws.Range("B64").Value = "Tag"
$ws.Range("C64").Value = "Haley Jones"
$ws.Range("D64").Value = "Rylee Mason"
$ws.Range("E64").Value = 45726
$ws.Range("F64").Value = 0.39861111111111114

# --- Row 65 : Action 64 -----------------------------------------------
$ws.Range("A65").Value = 64
$ws.Range("B65").Value = "Tag"
$ws.Range("C65").Value = "Alexis Pascual"
$ws.Range("D65").Value = "Colton Plank"
$ws.Range("E65").Value = 45726
$ws.Range("F65").Value = 0.39097222222222222

# --- Row 66 : trailing row, no Action # / Command ----------------------
$ws.Range("C66").Value = "Markus Laureano"
$ws.Range("D66").Value = "Damon Warwick"
$ws.Range("E66").Value = 45726
$ws.Range("F66").Value = 0.39374999999999999

# --- Unix-time formula, filled down (Excel groups these as one shared
#     formula, same as the existing G2:G28 / G29:G61 groups) -------------
$ws.Range("G62:G66").Formula = "=((E62+F62)-DATE(1970,1,1))*86400"

# --- Fix up formatting quirks so the new rows mirror how this sheet was
#     actually filled in (mix of fresh column-default formatting and a
#     few cells that were typed/copied with different/explicit formats) -

# A64/A65 continue the alternating s=5/s=7 "Action #" look used by A2:A61
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A64").PasteSpecial($xlPasteFormats)
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A65").PasteSpecial($xlPasteFormats)

# C65/C66 and D64/D66 were left as plain/general formatted cells (no
# border), matching cells such as C59/C60/D41/D58 earlier in the sheet
$ws.Range("D41").Copy() | Out-Null
$ws.Range("C65").PasteSpecial($xlPasteFormats)
$ws.Range("D41").Copy() | Out-Null
$ws.Range("C66").PasteSpecial($xlPasteFormats)
$ws.Range("D41").Copy() | Out-Null
$ws.Range("D64").PasteSpecial($xlPasteFormats)
$ws.Range("D41").Copy() | Out-Null
$ws.Range("D66").PasteSpecial($xlPasteFormats)

# G62:G66 keep the same explicit "Unix" number format/border as the rest
# of column G (the bare column-default style differs slightly in font)
$ws.Range("G61").Copy() | Out-Null
$ws.Range("G62:G66").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# --- Grow the table / autofilter to cover the new rows -----------------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:G66"))

# --- Update the on-screen selection/scroll position ---------------------
$ws.Range("A44").Select() | Out-Null
try {
    $excel.ActiveWindow.ScrollRow = 44
    $excel.ActiveWindow.ScrollColumn = 1
} catch {}
$ws.Range("H66").Select() | Out-Null
